$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '69.093.63'
$ws.Range("D3").Value = '3.744.79'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '3.742.75'
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.542'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("E10").Value = '  +3.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000249'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '4.368.26'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").Value = '3.758.55'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '69.105.43'
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.114'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.728'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("E24").Value = '  +8.71%  '
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("E31").Value = '  +1.53%  '
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("B34").Value = 'WrappedeETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D34").Value = '3.890.76'
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").Value = '3.677.33'
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E39").Value = '  +2.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.139'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.93%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '426.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("D50").Value = '2.787.13'
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("E51").Value = '  +0.15%  '
